$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 35. This shifts the existing rows 35..108
# down to 36..109 (all their cell content/styles move with them), exactly
# matching the rest of the diff (every row from the former 35 down to the
# former 108 reappears one row lower, and the sheet's dimension grows from
# R108 to R109).
$ws.Rows.Item(35).Insert()

# The newly inserted row 35 is blank; populate it with the new data record
# (same shape/columns as every other row in this table), using the new
# values for D/K/L/M/P called out in the diff and keeping the other
# columns identical to the rest of this "Haba" / Chillan series.
$ws.Cells.Item(35, 1).Value = 7
$ws.Cells.Item(35, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(35, 3).Value = "Ñuble"
$ws.Cells.Item(35, 4).Value = 45281
$ws.Cells.Item(35, 5).Value = 16
$ws.Cells.Item(35, 6).Value = 100112026
$ws.Cells.Item(35, 7).Value = "Haba"
$ws.Cells.Item(35, 8).Value = "Sin especificar"
$ws.Cells.Item(35, 9).Value = "Primera"
$ws.Cells.Item(35, 10).Value = 100
$ws.Cells.Item(35, 11).Value = 11000
$ws.Cells.Item(35, 12).Value = 12000
$ws.Cells.Item(35, 13).Value = 11500
$ws.Cells.Item(35, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(35, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(35, 16).Value = 460
$ws.Cells.Item(35, 17).Value = 25
$ws.Cells.Item(35, 18).Value = "Hortaliza"
